# Update produktion.xlsx til 14 dage: append a new day (2019-12-14, serial 43813)
# of hourly production data (rows 314-337) to Sheet1, extend the "data" defined
# name to cover the new rows, and move the selection to reflect where the user
# ended up after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(43813, 1, 468.7, 1393, 0, 2076.9),
    @(43813, 2, 451.6, 1375.9, 0, 2018.6),
    @(43813, 3, 451.2, 1433.5, 0, 1984.1),
    @(43813, 4, 432.5, 1511.3, 0, 1969.6),
    @(43813, 5, 427.8, 1625.5, 0, 1994.7),
    @(43813, 6, 446, 1646.7, 0, 2050),
    @(43813, 7, 510.8, 1479.8, 0, 2170.6999999999998),
    @(43813, 8, 539.9, 1597.9, 0, 2361),
    @(43813, 9, 567.1, 1671.4, 0, 2541.9),
    @(43813, 10, 580.5, 1596.8, 2, 2626.3),
    @(43813, 11, 599.6, 1412.5, 11, 2682),
    @(43813, 12, 607.1, 1283.2, 19, 2679.7),
    @(43813, 13, 606.29999999999995, 1254.2, 24, 2628.4),
    @(43813, 14, 592, 1342.8, 23, 2580.9),
    @(43813, 15, 576.20000000000005, 1487.5, 13, 2561.1),
    @(43813, 16, 570.5, 1776.7, 2, 2606.6),
    @(43813, 17, 576.79999999999995, 1843.6, 0, 2746.5),
    @(43813, 18, 563.70000000000005, 2032.7, 0, 2882.7),
    @(43813, 19, 567.4, 2324.8000000000002, 0, 2804.2),
    @(43813, 20, 539.70000000000005, 2509.4, 0, 2631),
    @(43813, 21, 445.3, 2677.9, 0, 2498.9),
    @(43813, 22, 321.8, 3023.8, 0, 2382.6),
    @(43813, 23, 293.60000000000002, 3221.4, 0, 2273.3000000000002),
    @(43813, 24, 256.2, 3345.9, 0, 2297.1999999999998)
)

$startRow = 314
$endRow = 337
$rows = $endRow - $startRow + 1
$cols = 6

# Copy the formatting (styles) of the last existing row down onto the new
# rows before writing values, so the new cells pick up the same date /
# hour / number styles as the rest of the table (s="3", s="4", s="5").
$ws.Range("A313:F313").Copy() | Out-Null
$ws.Range("A$startRow`:F$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Build a true 2-D array (PowerShell's @(@(...),@(...)) is a jagged array of
# objects, which Range.Value rejects) and assign it in one shot.
$arr = New-Object 'object[,]' $rows,$cols
for ($i = 0; $i -lt $rows; $i++) {
    for ($j = 0; $j -lt $cols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}
$ws.Range("A$startRow`:F$endRow").Value = $arr

# Extend the "data" named range to include the newly added rows.
$wb.Names("data").RefersTo = "=Sheet1!`$C`$2:`$F`$$endRow"

# Leave the selection where the author's cursor ended up after entering the
# last day's data.
$ws.Range("I334").Select() | Out-Null
